$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Shape 312" is the "deviceID/#" label box under the Server card.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Shape 312") {
        $shp = $s.Shapes.Item($i)
        break
    }
}

$tf  = $shp.TextFrame
$tr  = $tf.TextRange

# Split the text into "namespace/" + "deviceID" + "/#" across three runs,
# matching the authored edit ("shift to namespace first topics").
[void]$tr.InsertBefore("namespace/")
$full = $tf.TextRange
$deviceIdRun = $full.Characters(11, 8)
$deviceIdRun.Text = "deviceID"
$namespaceRun = $full.Characters(1, 10)
$namespaceRun.LanguageID = "en-US"

# The textbox uses wrap="none" + auto-fit, so widening the text also grows
# the shape; match the resulting position/size precisely.
$shp.Left   = 384.6111023622047
$shp.Top    = 176.69071197509768
$shp.Width  = 154.49417877197268
$shp.Height = 22.61874103546143
